$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.283.01"
$ws.Range("E2").Value = "  -3.60%  "

Set-TextValue $ws.Range("D3") "2.463.07"
$ws.Range("E3").Value = "  -2.73%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.11%  "

Set-TextValue $ws.Range("D5") "312.19"
$ws.Range("E5").Value = "  +0.44%  "

Set-TextValue $ws.Range("D6") "94.29"
$ws.Range("E6").Value = "  -6.51%  "

$ws.Range("E7").Value = "  -2.71%  "

$ws.Range("E8").Value = "  +0.05%  "

Set-TextValue $ws.Range("D9") "0.498"
$ws.Range("E9").Value = "  -4.90%  "

Set-TextValue $ws.Range("D10") "33.42"
$ws.Range("E10").Value = "  -6.47%  "

Set-TextValue $ws.Range("D11") "0.0779"
$ws.Range("E11").Value = "  -3.40%  "

$ws.Range("E12").Value = "  -1.09%  "

Set-TextValue $ws.Range("D13") "7.00"
$ws.Range("E13").Value = "  -4.83%  "

Set-TextValue $ws.Range("D14") "2.842.83"
$ws.Range("E14").Value = "  -2.64%  "

Set-TextValue $ws.Range("D15") "2.460.56"
$ws.Range("E15").Value = "  -1.97%  "

Set-TextValue $ws.Range("D16") "14.88"
$ws.Range("E16").Value = "  -3.44%  "

Set-TextValue $ws.Range("D17") "0.783"
$ws.Range("E17").Value = "  -4.12%  "

Set-TextValue $ws.Range("D18") "41.217.09"
$ws.Range("E18").Value = "  -3.71%  "

$ws.Range("E19").Value = "  -6.11%  "

$ws.Range("E20").Value = "  -3.35%  "

Set-TextValue $ws.Range("D21") "11.29"
$ws.Range("E21").Value = "  -8.53%  "

Set-TextValue $ws.Range("D22") "68.45"
$ws.Range("E22").Value = "  -1.54%  "

Set-TextValue $ws.Range("D23") "236.63"
$ws.Range("E23").Value = "  -2.95%  "

$ws.Range("E24").Value = "  -4.33%  "

$ws.Range("E25").Value = "  +0.05%  "

Set-TextValue $ws.Range("D26") "1.91"
$ws.Range("E26").Value = "  -6.37%  "

Set-TextValue $ws.Range("D27") "24.06"
$ws.Range("E27").Value = "  -5.65%  "

Set-TextValue $ws.Range("D28") "2.20"
$ws.Range("E28").Value = "  -6.10%  "

Set-TextValue $ws.Range("D29") "9.63"
$ws.Range("E29").Value = "  -5.44%  "

Set-TextValue $ws.Range("D30") "36.68"
$ws.Range("E30").Value = "  -5.68%  "

Set-TextValue $ws.Range("D31") "152.08"
$ws.Range("E31").Value = "  -5.61%  "

$ws.Range("E32").Value = "  -5.37%  "

$ws.Range("E33").Value = "  -5.11%  "

$ws.Range("E34").Value = "  -3.84%  "

Set-TextValue $ws.Range("D35") "0.0744"
$ws.Range("E35").Value = "  -5.73%  "

$ws.Range("E36").Value = "  -2.33%  "

Set-TextValue $ws.Range("B37") "ARBITRUM"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D37") "1.89"
$ws.Range("E37").Value = "  -4.34%  "

Set-TextValue $ws.Range("B38") "Celestia"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D38") "17.01"
$ws.Range("E38").Value = "  -7.86%  "

$ws.Range("E39").Value = "  -3.10%  "

$ws.Range("E40").Value = "  -7.73%  "

$ws.Range("E41").Value = "  +1.01%  "

$ws.Range("E42").Value = "  +0.19%  "

Set-TextValue $ws.Range("D43") "19.97"
$ws.Range("E43").Value = "  -10.86%  "

Set-TextValue $ws.Range("D44") "1.983.88"
$ws.Range("E44").Value = "  -0.48%  "

Set-TextValue $ws.Range("D45") "0.0284"
$ws.Range("E45").Value = "  -5.50%  "

Set-TextValue $ws.Range("D46") "3.04"
$ws.Range("E46").Value = "  -8.22%  "

Set-TextValue $ws.Range("D47") "8.81"
$ws.Range("E47").Value = "  -4.67%  "

Set-TextValue $ws.Range("B48") "Aave"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D48") "97.05"
$ws.Range("E48").Value = "  -3.87%  "

Set-TextValue $ws.Range("B49") "ordi"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue $ws.Range("D49") "68.95"
$ws.Range("E49").Value = "  -4.68%  "

Set-TextValue $ws.Range("D50") "0.178"
$ws.Range("E50").Value = "  -7.10%  "

Set-TextValue $ws.Range("D51") "74.38"
$ws.Range("E51").Value = "  -6.51%  "
